# Fixed plot saving bug
# Rows 202-401 (column A, the numeric index column) were off by one due to a
# previously deleted row whose gap was never closed. Decrement each of these
# rows' A value by 1 to close the gap (e.g. row202: 201->200 ... row401: 400->399).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 202; $r -le 401; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 1
}

# Row 242's title text also got its heart-emoji variation selector mangled
# into two Unicode replacement characters during the same save; reproduce
# that corrupted text exactly.
$heart = [char]0x2764
$fffd = [char]0xFFFD
$newTitle = "AMD Ryzen 5 5600X Desktop Processor, IN HAND" + $heart + $fffd + $fffd + ", WHO WILL BE LUCKY?!?"
$ws.Cells.Item(242, 2).Value2 = $newTitle
